$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting of the last fully-populated data row (row 32) down into
# rows 33-36 so the new rows pick up the same "data row" styles (date format,
# bottom-border, etc.) instead of the blank "filler" row styles.
$ws.Range("C32:L32").Copy()
$ws.Range("C33:C36").PasteSpecial(-4122) | Out-Null
$ws.Range("D33:D36").PasteSpecial(-4122) | Out-Null
$ws.Range("E33:E36").PasteSpecial(-4122) | Out-Null
$ws.Range("G33:G36").PasteSpecial(-4122) | Out-Null
$ws.Range("I33:I36").PasteSpecial(-4122) | Out-Null
$ws.Range("K33:K36").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Cambridge 13 Test 1-4: four new practiced tests entered into the log.
$ws.Range("C33").Value = 30
$ws.Range("D33").Value = 45512
$ws.Range("E33").Value = "IELTS13_Test1"
$ws.Range("J33").Value = 1.1

$ws.Range("C34").Value = 31
$ws.Range("D34").Value = 45513
$ws.Range("E34").Value = "IELTS13_Test2"
$ws.Range("J34").Value = 1.1

$ws.Range("C35").Value = 32
$ws.Range("D35").Value = 45514
$ws.Range("E35").Value = "IELTS13_Test3"
$ws.Range("J35").Value = 1.1

$ws.Range("C36").Value = 33
$ws.Range("D36").Value = 45515
$ws.Range("E36").Value = "IELTS13_Test4"

$ws.Range("J36").Select()
